$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the TC1 - step 2 (row 11) to the new "ordenar pelo nome do servidor" step.
#    (Do this BEFORE duplicating the block, so the duplicate below doesn't
#    inherit the old wording.)
$ws.Range("B11").Value = "Chefe Clica para ordenar pelo nome do servidor."
$ws.Range("D11").Value = "SYSTEM Visualiza os registros de solicitações de diária ordenado pelo nome do servidor."

# 2. Insert 8 fresh rows right after TC1's block (before the first blank
#    separator row 12) to make room for a brand-new TC2 block. This shifts
#    every row from 12 onward down by 8 (mirroring the diff's row shift),
#    carrying along styles and merged cells automatically.
$ws.Rows("12:19").Insert()

# 3. Populate the new TC2 block (rows 14-19) by duplicating the TC1 block
#    structure (rows 6-11), which has identical formatting/merges, then fix
#    up the handful of cells that differ for TC2.
$ws.Range("A6:F11").Copy($ws.Range("A14"))

# Rows 12-13 must stay as plain blank separator rows (no inherited style).
$ws.Range("A12:F13").ClearContents()
$ws.Range("A12:F13").ClearFormats()

# TC2 header and step text (this is the step that used to live on TC1).
$ws.Range("B14").Value = "TC2"
$ws.Range("B19").Value = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$ws.Range("D19").Value = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."

# 4. Relabel the old TC2/TC3/TC4 blocks, which shifted down by 8 rows and
#    now need to read TC3/TC4/TC5 respectively.
$ws.Range("B22").Value = "TC3"
$ws.Range("B30").Value = "TC4"
$ws.Range("B38").Value = "TC5"

# 5. Bump the reported test-suite size from 4 to 5 test cases.
$ws.Range("D3").Value = "Size: 5 test case(s))"
